# Apply scheduled-runner value updates to the per-job-sheet profit tables
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR), matching the Alpha_Profits.xlsx OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 964.1
$ws.Range("I15").Value = 964.1
$ws.Range("K15").Value = 2892.3
$ws.Range("M15").Value = -2723.3
# Row 100
$ws.Range("H100").Value = 4239.6
$ws.Range("I100").Value = 2549.5
$ws.Range("K100").Value = 2549.5
$ws.Range("M100").Value = -2008.5
# Row 107
$ws.Range("H107").Value = 762.8570999999999
$ws.Range("I107").Value = 762.8570999999999
$ws.Range("K107").Value = 762.8570999999999
$ws.Range("M107").Value = 1157.1429
# Row 132
$ws.Range("H132").Value = 63961.125
$ws.Range("I132").Value = 68090.07000000001
$ws.Range("J132").Value = 2027
$ws.Range("K132").Value = 204270.21
$ws.Range("L132").Value = 6081
$ws.Range("M132").Value = -201740.21
$ws.Range("N132").Value = -11141
# Row 137
$ws.Range("H137").Value = 6484.6924
$ws.Range("I137").Value = 2866.4443
$ws.Range("K137").Value = 8599.332900000001
$ws.Range("M137").Value = -6049.332900000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1232.4857
$ws.Range("I32").Value = 1239.3235
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 1239.3235
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -952.3235
$ws.Range("N32").Value = -1574
# Row 45
$ws.Range("H45").Value = 1122.5454
$ws.Range("I45").Value = 994.8
$ws.Range("K45").Value = 994.8
$ws.Range("M45").Value = -617.8
# Row 74
$ws.Range("H74").Value = 6615475
$ws.Range("I74").Value = 3705509.8
$ws.Range("J74").Value = 13890388
$ws.Range("K74").Value = 3705509.8
$ws.Range("L74").Value = 13890388
$ws.Range("M74").Value = -3704635.8
$ws.Range("N74").Value = -13892136
# Row 77
$ws.Range("H77").Value = 6615475
$ws.Range("I77").Value = 3705509.8
$ws.Range("J77").Value = 13890388
$ws.Range("K77").Value = 18527549
$ws.Range("L77").Value = 69451940
$ws.Range("M77").Value = -18523181
$ws.Range("N77").Value = -69460676
# Row 102
$ws.Range("H102").Value = 822.5
$ws.Range("I102").Value = 629
$ws.Range("K102").Value = 629
$ws.Range("M102").Value = 993
# Row 122
$ws.Range("H122").Value = 1792.5264
$ws.Range("I122").Value = 1607.7273
$ws.Range("J122").Value = 2046.625
$ws.Range("K122").Value = 4823.1819
$ws.Range("L122").Value = 6139.875
$ws.Range("M122").Value = -2373.1819
$ws.Range("N122").Value = -11039.875
# Row 132
$ws.Range("H132").Value = 38464190
$ws.Range("I132").Value = 2831.889
$ws.Range("J132").Value = 125002250
$ws.Range("K132").Value = 8495.667000000001
$ws.Range("L132").Value = 375006750
$ws.Range("M132").Value = -5965.667000000001
$ws.Range("N132").Value = -375011810

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 936.25
$ws.Range("I20").Value = 936.25
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 936.25
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -689.25
$ws.Range("N20").ClearContents()
# Row 86
$ws.Range("H86").Value = 2789.6428
$ws.Range("I86").Value = 2274.4
$ws.Range("J86").Value = 3075.889
$ws.Range("K86").Value = 2274.4
$ws.Range("L86").Value = 3075.889
$ws.Range("M86").Value = -1151.4
$ws.Range("N86").Value = -5321.889
# Row 89
$ws.Range("H89").Value = 2789.6428
$ws.Range("I89").Value = 2274.4
$ws.Range("J89").Value = 3075.889
$ws.Range("K89").Value = 11372
$ws.Range("L89").Value = 15379.445
$ws.Range("M89").Value = -5756
$ws.Range("N89").Value = -26611.445
# Row 94
$ws.Range("H94").Value = 2053.4119
$ws.Range("I94").Value = 1666.5333
$ws.Range("K94").Value = 1666.5333
$ws.Range("M94").Value = -1215.5333
# Row 99
$ws.Range("H99").Value = 2392.8572
$ws.Range("J99").Value = 1940.6666
$ws.Range("L99").Value = 1940.6666
$ws.Range("N99").Value = -4936.6666
# Row 134
$ws.Range("H134").Value = 23335522
$ws.Range("I134").Value = 11906892
$ws.Range("K134").Value = 35720676
$ws.Range("M134").Value = -35718141

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2075.2942
$ws.Range("I31").Value = 2006.6428
$ws.Range("K31").Value = 2006.6428
$ws.Range("M31").Value = -1711.6428
# Row 34
$ws.Range("H34").Value = 2075.2942
$ws.Range("I34").Value = 2006.6428
$ws.Range("K34").Value = 2006.6428
$ws.Range("M34").Value = -1804.6428
# Row 58
$ws.Range("H58").Value = 1361.069
$ws.Range("I58").Value = 1361.069
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1361.069
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1158.069
$ws.Range("N58").ClearContents()
# Row 94
$ws.Range("H94").Value = 1429.1666
$ws.Range("I94").Value = 1100
$ws.Range("K94").Value = 1100
$ws.Range("M94").Value = -649
# Row 132
$ws.Range("H132").Value = 6952.72
$ws.Range("I132").Value = 7071.85
$ws.Range("J132").Value = 6476.2
$ws.Range("K132").Value = 21215.55
$ws.Range("L132").Value = 19428.6
$ws.Range("M132").Value = -18685.55
$ws.Range("N132").Value = -24488.6
# Row 136
$ws.Range("H136").Value = 1361.069
$ws.Range("I136").Value = 1361.069
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4083.207
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1533.207
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 935.1667
$ws.Range("I8").Value = 935.1667
$ws.Range("K8").Value = 2805.5001
$ws.Range("M8").Value = -2666.5001

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2097
$ws.Range("I102").Value = 2119.8462
$ws.Range("K102").Value = 2119.8462
$ws.Range("M102").Value = -497.8462

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2282.1428
$ws.Range("I7").Value = 1996.6666
$ws.Range("K7").Value = 1996.6666
$ws.Range("M7").Value = -1884.6666
# Row 100
$ws.Range("H100").Value = 1703
$ws.Range("I100").Value = 1703
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1703
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1162
$ws.Range("N100").ClearContents()
# Row 122
$ws.Range("H122").Value = 3547.923
$ws.Range("I122").Value = 3312.6
$ws.Range("K122").Value = 9937.799999999999
$ws.Range("M122").Value = -7487.799999999999
# Row 126
$ws.Range("H126").Value = 2282.1428
$ws.Range("I126").Value = 1996.6666
$ws.Range("K126").Value = 5989.9998
$ws.Range("M126").Value = -3519.9998

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 20362.572
$ws.Range("I54").Value = 20076
$ws.Range("J54").Value = 20410.334
$ws.Range("K54").Value = 20076
$ws.Range("L54").Value = 20410.334
$ws.Range("M54").Value = -19556
$ws.Range("N54").Value = -21450.334
# Row 103
$ws.Range("H103").Value = 37500
$ws.Range("J103").Value = 37500
$ws.Range("L103").Value = 37500
$ws.Range("N103").Value = -39844
